$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column B (stress level responses) for rows 2-7
$ws.Range("B2").Value = "A little stressful "
$ws.Range("B3").Value = "Very stressful "
$ws.Range("B4").Value = "Moderately stressful "
$ws.Range("B5").Value = "Very stressful "
$ws.Range("B6").Value = "Very stressful "
$ws.Range("B7").Value = "Very stressful "

# Update the active selection to match the author's final cursor position
$ws.Range("B8").Select()
